$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 2-10: new reporting date (2024-08-10 -> 2024-08-13) and refreshed station load figures ---

# Force column A to stay textual so the date string is not auto-converted to a date serial number
$ws.Range("A2:A10").NumberFormat = "@"

$ws.Range("A2").Value = "2024-08-13"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 22060
$ws.Range("D2").Value = 6063
$ws.Range("E2").Value = 9842
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 5734
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 17979
$ws.Range("J2").Value = 16944
$ws.Range("K2").Value = 78622
$ws.Range("L2").Value = 22062.206
$ws.Range("M2").Value = 6063.6063
$ws.Range("N2").Value = 9842.984200000001
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 5734.5734
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 17980.7979
$ws.Range("S2").Value = 16945.6944
$ws.Range("T2").Value = 78629.8622

$ws.Range("A3").Value = "2024-08-13"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 21151
$ws.Range("D3").Value = 5712
$ws.Range("E3").Value = 9409
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 5549
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 17214
$ws.Range("J3").Value = 16164
$ws.Range("K3").Value = 75199
$ws.Range("L3").Value = 21153.1151
$ws.Range("M3").Value = 5712.5712
$ws.Range("N3").Value = 9409.9409
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 5549.5549
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 17215.7214
$ws.Range("S3").Value = 16165.6164
$ws.Range("T3").Value = 75206.5199

$ws.Range("A4").Value = "2024-08-13"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 20342
$ws.Range("D4").Value = 5485
$ws.Range("E4").Value = 9148
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 5363
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 14058
$ws.Range("J4").Value = 15424
$ws.Range("K4").Value = 69820
$ws.Range("L4").Value = 20344.0342
$ws.Range("M4").Value = 5485.5485
$ws.Range("N4").Value = 9148.9148
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 5363.5363
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 14059.4058
$ws.Range("S4").Value = 15425.5424
$ws.Range("T4").Value = 69826.982

$ws.Range("A5").Value = "2024-08-13"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 19725
$ws.Range("D5").Value = 5175
$ws.Range("E5").Value = 8877
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 5217
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 15883
$ws.Range("J5").Value = 14832
$ws.Range("K5").Value = 69709
$ws.Range("L5").Value = 19726.9725
$ws.Range("M5").Value = 5175.5175
$ws.Range("N5").Value = 8877.887699999999
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 5217.5217
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 15884.5883
$ws.Range("S5").Value = 14833.4832
$ws.Range("T5").Value = 69715.9709

$ws.Range("A6").Value = "2024-08-13"
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 19604
$ws.Range("D6").Value = 5284
$ws.Range("E6").Value = 8705
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 5424
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 15163
$ws.Range("J6").Value = 14968
$ws.Range("K6").Value = 69148
$ws.Range("L6").Value = 19605.9604
$ws.Range("M6").Value = 5284.5284
$ws.Range("N6").Value = 8705.870499999999
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 5424.5424
$ws.Range("Q6").Value = 0
$ws.Range("R6").Value = 15164.5163
$ws.Range("S6").Value = 14969.4968
$ws.Range("T6").Value = 69154.91479999998

$ws.Range("A7").Value = "2024-08-13"
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 19818
$ws.Range("D7").Value = 5342
$ws.Range("E7").Value = 9157
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 5658
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 15700
$ws.Range("J7").Value = 15351
$ws.Range("K7").Value = 71026
$ws.Range("L7").Value = 19819.9818
$ws.Range("M7").Value = 5342.5342
$ws.Range("N7").Value = 9157.9157
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 5658.5658
$ws.Range("Q7").Value = 0
$ws.Range("R7").Value = 15701.57
$ws.Range("S7").Value = 15352.5351
$ws.Range("T7").Value = 71033.1026

$ws.Range("A8").Value = "2024-08-13"
$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 19108
$ws.Range("D8").Value = 4547
$ws.Range("E8").Value = 9238
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 5986
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 15087
$ws.Range("J8").Value = 15403
$ws.Range("K8").Value = 69369
$ws.Range("L8").Value = 19109.9108
$ws.Range("M8").Value = 4547.4547
$ws.Range("N8").Value = 9238.9238
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 5986.5986
$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = 15088.5087
$ws.Range("S8").Value = 15404.5403
$ws.Range("T8").Value = 69375.9369

$ws.Range("A9").Value = "2024-08-13"
$ws.Range("B9").Value = 8
$ws.Range("C9").Value = 22031
$ws.Range("D9").Value = 4781
$ws.Range("E9").Value = 10950
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 9539
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 18194
$ws.Range("J9").Value = 17943
$ws.Range("K9").Value = 83438
$ws.Range("L9").Value = 22033.2031
$ws.Range("M9").Value = 4781.4781
$ws.Range("N9").Value = 10951.095
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 9539.9539
$ws.Range("Q9").Value = 0
$ws.Range("R9").Value = 18195.8194
$ws.Range("S9").Value = 17944.7943
$ws.Range("T9").Value = 83446.3438

$ws.Range("A10").Value = "2024-08-13"
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 28392
$ws.Range("D10").Value = 5660
$ws.Range("E10").Value = 13268
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 20628
$ws.Range("K10").Value = 67948
$ws.Range("L10").Value = 28394.8392
$ws.Range("M10").Value = 5660.566
$ws.Range("N10").Value = 13269.3268
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 20630.0628
$ws.Range("T10").Value = 67954.7948

# Remove the leftover default-format marker on column A now that the text values are committed
$ws.Range("A2:A10").ClearFormats()

# --- Remove rows 11 and 12 (hour 10 and hour 11); the table now only spans hours 1-9 (A1:T10) ---
$ws.Rows("11:12").Delete()

